$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96 (pushes existing rows 96-208 down to 97-209)
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new record
$ws.Cells.Item(96, 1).Value = 11
$ws.Cells.Item(96, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(96, 3).Value = "Bíobío"
$ws.Cells.Item(96, 4).Value = 44579
$ws.Cells.Item(96, 5).Value = 8
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100102
$ws.Cells.Item(96, 8).Value = "Cítricos"
$ws.Cells.Item(96, 9).Value = 100102005
$ws.Cells.Item(96, 10).Value = "Naranja"
$ws.Cells.Item(96, 11).Value = "Valencia"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 270
$ws.Cells.Item(96, 14).Value = 9000
$ws.Cells.Item(96, 15).Value = 9500
$ws.Cells.Item(96, 16).Value = 9222
$ws.Cells.Item(96, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(96, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(96, 19).Value = 615
$ws.Cells.Item(96, 20).Value = 15
